# Generate Report for Handback
#
# Source (en-US) handback files got new GUID-based names:
#   67366a8b-6080-41fa-b81b-6cb6d38cd198  ->  90e85150-993a-4084-ae5c-840ad9e24d3a
#   c1417f25-d57e-4b55-9a9d-02eb1f70efd5  ->  fffff3d00647-7864-490d-b741-2343d10cff78
# and their corresponding handoff/handback package (.xlf) content hash changed to:
#   4c1798fc505cf4fb24658657e17b4eaffe53105d / 5f77132c4b14e7c9257815f69be52939f817dab4
#     -> d808a6c23d42fdb3ce1773c39e77884d0bfc06f3
# (note: both rows now resolve to the SAME .xlf package name/hash)
# along with refreshed handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "67366a8b-6080-41fa-b81b-6cb6d38cd198"
$newGuid1 = "90e85150-993a-4084-ae5c-840ad9e24d3a"
$oldGuid2 = "c1417f25-d57e-4b55-9a9d-02eb1f70efd5"
$newGuid2 = "fffff3d00647-7864-490d-b741-2343d10cff78"
$newHash  = "d808a6c23d42fdb3ce1773c39e77884d0bfc06f3"

# ---------------------------------------------------------------------------
# Sheet "Overview": just the two .md file-name cells
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("A3").Value = "$newGuid2.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 (90e85150...)
$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("D2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-21 21:03:59"
$ws.Range("F2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-21 21:04:20"

# Row 3 (fffff3d00647... but the handoff/handback package now matches row 2's)
$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("D3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-21 21:03:59"
$ws.Range("F3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-21 21:04:20"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 (90e85150...)
$ws.Range("A2").Value = "$newGuid1.md"
$ws.Range("D2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("E2").Value = "2016-03-21 21:04:03"
$ws.Range("F2").Value = "$newGuid1.md"
$ws.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-03-21 21:04:27"

# Row 3 (fffff3d00647... but the handoff/handback package now matches row 2's)
$ws.Range("A3").Value = "$newGuid2.md"
$ws.Range("D3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("E3").Value = "2016-03-21 21:04:03"
$ws.Range("F3").Value = "$newGuid2.md"
$ws.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$ws.Range("H3").Value = "2016-03-21 21:04:27"

# ---------------------------------------------------------------------------
# Hyperlink display text: keep the underlying link targets untouched (the
# diff does not touch the .rels files) but refresh the cached display text
# shown for each linked cell so it matches the new file names.
# ---------------------------------------------------------------------------
function Set-HyperlinkDisplay($sheetName, $cellRef, $text) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $range = $sheet.Range($cellRef)
    if ($range.Hyperlinks.Count -gt 0) {
        $range.Hyperlinks.Item(1).TextToDisplay = $text
    }
}

Set-HyperlinkDisplay "Overview" "A2" "$newGuid1.md"
Set-HyperlinkDisplay "Overview" "A3" "$newGuid2.md"

Set-HyperlinkDisplay "zh-cn" "A2" "$newGuid1.md"
Set-HyperlinkDisplay "zh-cn" "D2" "$newGuid1.$newHash.zh-cn.xlf"
Set-HyperlinkDisplay "zh-cn" "F2" "$newGuid1.md"
Set-HyperlinkDisplay "zh-cn" "G2" "$newGuid1.$newHash.zh-cn.xlf"
Set-HyperlinkDisplay "zh-cn" "A3" "$newGuid2.md"
Set-HyperlinkDisplay "zh-cn" "D3" "$newGuid1.$newHash.zh-cn.xlf"
Set-HyperlinkDisplay "zh-cn" "F3" "$newGuid2.md"
Set-HyperlinkDisplay "zh-cn" "G3" "$newGuid1.$newHash.zh-cn.xlf"

Set-HyperlinkDisplay "de-de" "A2" "$newGuid1.md"
Set-HyperlinkDisplay "de-de" "D2" "$newGuid1.$newHash.de-de.xlf"
Set-HyperlinkDisplay "de-de" "F2" "$newGuid1.md"
Set-HyperlinkDisplay "de-de" "G2" "$newGuid1.$newHash.de-de.xlf"
Set-HyperlinkDisplay "de-de" "A3" "$newGuid2.md"
Set-HyperlinkDisplay "de-de" "D3" "$newGuid1.$newHash.de-de.xlf"
Set-HyperlinkDisplay "de-de" "F3" "$newGuid2.md"
Set-HyperlinkDisplay "de-de" "G3" "$newGuid1.$newHash.de-de.xlf"
